# NIT-9015643408.xlsx — "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"
#
# The "Estado de Cuenta" table (rows 16-21) lists one debt period per row for
# worker PPT/5964707 WILLIAM JESUS ROJAS (periods 2503..2508). This edit adds
# a new period (2509) as a new last row of the table, and refreshes the
# summary totals above the table ("VALOR MORA" and "Cant. Periodos").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new table row (row 22) for period 2509 --------------------------
# Insert a blank row right after the current last data row (21); this shifts
# the trailing "firma" rows (26/27) down to (27/28), same as the target.
$ws.Rows.Item(22).Insert()

# The new row 22 becomes the table's new last row, so it should carry the
# "closing" (bottom-border) formatting that row 21 used to have. Copying
# row 21 (values + format) into the fresh row 22 achieves both at once.
$ws.Range("B21:J21").Copy($ws.Range("B22:J22"))

# Row 21 is no longer the last row of the table, so it reverts to the same
# "interior" formatting used by the other body rows (16-20). Copy format
# from row 20 to pick that up.
$ws.Range("B20:J20").Copy($ws.Range("B21:J21"))

# The copy from row 20 overwrote row 21's period text with row 20's value
# ("2507") — put row 21's own period code back.
$ws.Range("E21").Value = "2508"

# Fill in the new row's period code (everything else — Tipo Doc, N° Doc,
# Nombre, Valor Mora, Salario Basico — is identical across all periods, and
# was already copied from row 21 above).
$ws.Range("E22").Value = "2509"

# --- Refresh the summary figures above the table ----------------------------
# VALOR MORA (total overdue amount) grows from 341640 to 398580 with the new
# period's debt.
$ws.Range("E11").Value = 398580

# Cant. Periodos (count of overdue periods) goes from 6 to 7.
$ws.Range("F13").Value = 7

Write-Host "Added period 2509 row; updated VALOR MORA and Cant. Periodos totals."
